$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "181÷8=22, 5" "109÷7=15, 4"
Replace-Text "518÷7=74, 0" "323÷2=161, 1"
Replace-Text "196÷2=98, 0" "186÷8=23, 2"
Replace-Text "985÷8=123, 1" "158÷7=22, 4"
Replace-Text "557÷6=92, 5" "902÷2=451, 0"
Replace-Text "112÷9=12, 4" "577÷5=115, 2"
Replace-Text "126÷3=42, 0" "827÷6=137, 5"
Replace-Text "706÷7=100, 6" "931÷9=103, 4"
Replace-Text "653÷9=72, 5" "501÷8=62, 5"
Replace-Text "884÷9=98, 2" "230÷3=76, 2"
Replace-Text "174÷6=29, 0" "363÷2=181, 1"
Replace-Text "837÷9=93, 0" "427÷2=213, 1"
Replace-Text "985÷4=246, 1" "101÷6=16, 5"
Replace-Text "111÷3=37, 0" "928÷6=154, 4"
Replace-Text "292÷4=73, 0" "296÷9=32, 8"
Replace-Text "250÷4=62, 2" "669÷4=167, 1"
Replace-Text "724÷4=181, 0" "607÷8=75, 7"
Replace-Text "505÷7=72, 1" "320÷3=106, 2"
Replace-Text "174÷2=87, 0" "701÷8=87, 5"
Replace-Text "782÷8=97, 6" "674÷6=112, 2"
Replace-Text "628÷8=78, 4" "931÷7=133, 0"
Replace-Text "426÷3=142, 0" "442÷3=147, 1"
Replace-Text "532÷2=266, 0" "325÷2=162, 1"
Replace-Text "841÷2=420, 1" "998÷8=124, 6"
Replace-Text "358÷4=89, 2" "994÷3=331, 1"

Write-Output "Done"
